# Apply the "implemented excel import experimentally" changes:
#  - Add a new Material row (Mineral wool / Insulation) to the "Materials" sheet
#  - Add two new Construction rows ("Construction 1" and "Construction 2") to
#    the "Constructions" sheet, referencing the new "Mineral wool" material
#  - Update the active selections on both sheets
#  - Make "Constructions" the active/selected sheet (it was "Materials" before)

$wb = $excel.ActiveWorkbook

$wsMaterials = $wb.Worksheets.Item("Materials")
$wsConstructions = $wb.Worksheets.Item("Constructions")

# --- Materials sheet: new row 2 -------------------------------------------------
$wsMaterials.Range("A2").Value = "Mineral wool"
$wsMaterials.Range("B2").Value = 0.7
$wsMaterials.Range("C2").Value = 1030
$wsMaterials.Range("D2").Value = 0.7
$wsMaterials.Range("E2").Value = "Insulation"
$wsMaterials.Range("F2").Value = 0.9
$wsMaterials.Range("G2").Value = 100
$wsMaterials.Range("H2").Value = 0.036
$wsMaterials.Range("I2").Value = 0.35
$wsMaterials.Range("J2").Value = 0.05
$wsMaterials.Range("K2").Value = 0.02

# --- Constructions sheet: new rows 2 & 3 ----------------------------------------
$wsConstructions.Range("A2").Value = "Construction 1"
$wsConstructions.Range("B2").Value = "Mineral wool"
$wsConstructions.Range("C2").Value = 0.3

$wsConstructions.Range("A3").Value = "Construction 2"
$wsConstructions.Range("B3").Value = "Mineral wool"
$wsConstructions.Range("C3").Value = 0.3
$wsConstructions.Range("D3").Value = "Mineral wool"
$wsConstructions.Range("E3").Value = 0.1

# --- Update selections on each sheet --------------------------------------------
$null = $wsMaterials.Range("I3").Select()
$null = $wsConstructions.Range("A4").Select()

# --- Make "Constructions" the active sheet (tabSelected moves to it) -----------
$null = $wsConstructions.Activate()
